$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for data rows 2-45
# from 2023-10-09 (serial 45208) to 2023-10-13 (serial 45212)
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
